# Apply model def file structure updates:
# On the "comp_type_dmg_algo" sheet, insert two new columns:
#  - "location" (after "median", before "beta") with value 0 for every data row
#  - "recovery_function" (after "lower_limit", before "recovery_mean") with value "Normal" for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comp_type_dmg_algo")

# Find the last used row (29 in the original layout: header + 28 data rows)
$lastRow = $ws.Cells(1, 1).EntireColumn.Cells.SpecialCells(11).Row

# --- Insert "recovery_function" column before current column M (recovery_mean) ---
$ws.Columns("M").Insert(-4161) | Out-Null
$ws.Range("M1").Value = "recovery_function"
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 13).Value = "Normal"  # column M -> recovery_function
}

# --- Insert "location" column before current column G (beta) ---
$ws.Columns("G").Insert(-4161) | Out-Null
$ws.Range("G1").Value = "location"
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 7).Value = 0          # column G -> location
}

$wb.Save()
